$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.185.24'
$ws.Range('D3').Value = '3.870.11'
$ws.Range('E3').Value = '  +1.60%  '
$ws.Range('E4').Value = '  +0.10%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '698.58'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.14%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '174.19'
$c.Style = 'Normal'
$ws.Range('D7').Value = '3.867.35'
$ws.Range('E7').Value = '  +1.55%  '
$ws.Range('E8').Value = '  +0.03%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.526'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('E11').Value = '  -8.29%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.460'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.50%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.0000262'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +3.16%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '36.53'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +0.64%  '
$ws.Range('D15').Value = '4.522.49'
$ws.Range('E15').Value = '  +1.61%  '
$ws.Range('D16').Value = '3.874.04'
$ws.Range('E16').Value = '  +1.60%  '
$ws.Range('D17').Value = '71.241.92'
$ws.Range('E17').Value = '  +0.45%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '17.70'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.38%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '7.25'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.50%  '
$ws.Range('E20').Value = '  -0.33%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '11.17'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -1.57%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '499.86'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +4.05%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.724'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +1.10%  '
$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.0000150'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +2.09%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '85.01'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +4.49%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '12.27'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.97%  '
$ws.Range('E28').Value = '  +0.18%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '3.19'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +2.73%  '
$ws.Range('E30').Value = '  -0.02%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '7.59'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +1.23%  '
$ws.Range('E32').Value = '  -1.46%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '29.78'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.49%  '
$ws.Range('E34').Value = '  +1.53%  '
$ws.Range('E35').Value = '  +0.70%  '
$ws.Range('D36').Value = '3.821.90'
$ws.Range('E36').Value = '  +1.70%  '
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('E38').Value = '  +2.43%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '2.39'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +8.68%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '3.45'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -3.29%  '
$ws.Range('E41').Value = '  +8.39%  '
$ws.Range('E42').Value = '  +0.67%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('E45').Value = '  -7.63%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '163.73'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +2.23%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '49.08'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.54%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '419.14'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +4.42%  '
$ws.Range('B49').Value = 'TheGraph'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.303'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +1.31%  '
$ws.Range('B50').Value = 'Arweave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '43.78'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -3.06%  '
$ws.Range('B51').Value = 'ONDO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '1.39'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -2.56%  '
